$wb = $excel.ActiveWorkbook

$wsRoles = $wb.Worksheets.Item("Roles&Rules")
$wsAlex  = $wb.Worksheets.Item("Alex")

# --- Alex sheet: fill in desirable-roles / requirements table ---------------
# Copy the boxed-border style already used on the Roles&Rules sheet (B2) onto
# the three role cells B2:B4 before writing their text.
$wsRoles.Range("B2").Copy()
$wsAlex.Range("B2:B4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsAlex.Cells.Item(2, 2).Value = "Tech Researcher"
$wsAlex.Cells.Item(2, 3).Value = "Critical thinking, programming"

$wsAlex.Cells.Item(3, 2).Value = "Coordinator"
$wsAlex.Cells.Item(3, 3).Value = "Time management"

$wsAlex.Cells.Item(4, 2).Value = "Editor"
$wsAlex.Cells.Item(4, 3).Value = "Attention to detail, critical thinking"

$wsAlex.Cells.Item(5, 2).Value = "Programming, attention to detail"

$wsAlex.Cells.Item(8, 2).Value = "Computer science, machine learning projects, other computer science projects"

$wsAlex.Cells.Item(10, 1).Value = "Add/change"
$wsAlex.Cells.Item(10, 2).Value = "I think there should be a role dedicated to creating the solution and relying on the tech researcher's research"

# --- Selections / active sheet ----------------------------------------------
$wsRoles.Range("C6").Select()
$wsAlex.Range("E13").Select()

# Make "Alex" the active (visible / selected) sheet/tab.
$wsAlex.Activate()
